$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the unit price column (D28:D34) with new computed values
$ws.Range("D28").Value = 636.951
$ws.Range("D29").Value = 667.362
$ws.Range("D30").Value = 1075.196
$ws.Range("D31").Value = 1526.765
$ws.Range("D32").Value = 2133.272
$ws.Range("D33").Value = 2737.882
$ws.Range("D34").Value = 4277.943
